$d = $word.ActiveDocument

# --- Paragraph 1: replace the placeholder methods text -------------------
$d.Content.Find.Execute(
    "TO DO. describe climdb methods in a few paragraphs.",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "When the database was active, contributors interacting with the ClimHydroDB database referred to the User Guide (PDF) archived in this package for methods and required steps. The User Guide also includes some details about required fields and flagging.",
    2
)

# --- Paragraph 2: was an empty paragraph, now gets new text ---------------
$p2 = $d.Paragraphs.Item(2)
$p2.Range.Text = "Historical documents, earlier versions of user guides, and presentations describing the database can be found in the zip of related material, also in this archive."

# --- Paragraph 3: replace "This is paragraph 2." text ---------------------
$d.Content.Find.Execute(
    "This is paragraph 2.",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "The export of data and creation of this archive was accomplished as part of EDI and LTER’s Clim-HydroDB-2.0 project. Scripts, raw metadata and related material can be found in the GitHub repository: https://github.com/lter/Clim-HydroDB-2.0. ",
    2
)

# --- Remove the old paragraph 4 (empty) and paragraph 5 ("This is paragraph 3.") ---
# After the text replacement above, paragraph indices are unchanged (replace
# was in-place). The empty paragraph 4 and the now-orphaned "paragraph 3"
# text paragraph (index 5) are no longer needed and get merged away.
$d.Paragraphs.Item(4).Range.Delete()
$d.Paragraphs.Item(4).Range.Delete()

# --- Apply paragraph spacing (w:after=200, w:lineRule=auto) to the three
#     surviving content paragraphs --------------------------------------
for ($i = 1; $i -le 3; $i++) {
    $p = $d.Paragraphs.Item($i)
    $p.SpaceAfter = 10
    $p.LineSpacingRule = 5
}
